# Update "想去人数" (column F) counters on several rows across sheets,
# matching the content refresh recorded in commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 616
$ws.Range("F6").Value = 6139
$ws.Range("F7").Value = 680
$ws.Range("F9").Value = 35
$ws.Range("F13").Value = 617
$ws.Range("F15").Value = 64
$ws.Range("F17").Value = 308
$ws.Range("F18").Value = 1385
$ws.Range("F27").Value = 3409

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F14").Value = 89
$ws.Range("F19").Value = 367

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F8").Value = 1529
$ws.Range("F12").Value = 699

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value = 1529
$ws.Range("F9").Value = 699
$ws.Range("F10").Value = 616
$ws.Range("F14").Value = 6139
$ws.Range("F16").Value = 680
$ws.Range("F18").Value = 35
$ws.Range("F22").Value = 617
$ws.Range("F25").Value = 89
$ws.Range("F29").Value = 64
$ws.Range("F31").Value = 308
$ws.Range("F33").Value = 367
$ws.Range("F34").Value = 1385
$ws.Range("F48").Value = 3409
